$d = $word.ActiveDocument

function Replace-InCell($tableIndex, $row, $col, $findText, $replaceText) {
    $table = $d.Tables.Item($tableIndex)
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    # Exclude the trailing end-of-cell marker from the range so the
    # replacement stays confined to the cell's own text.
    $scoped = $d.Range($cellRange.Start, $cellRange.End - 1)
    # wdReplaceOne (1) keeps the operation limited to the supplied range;
    # wdReplaceAll would replace every match in the whole document.
    $scoped.Find.ClearFormatting()
    $scoped.Find.Execute($findText, $true, $false, $false, $false, $false, `
                          $true, 1, $false, $replaceText, 1)
}

# "F" + "K" runs in the header row of the KORISNIK table merge into one "FK" run.
Replace-InCell 2 1 5 "FK" "FK"

# kor_email: nvarchar(256) -> nvarchar(128)
Replace-InCell 2 4 3 "nvarchar(256)" "nvarchar(128)"

# kor_pwdhash: binary(128) -> binary(16)
Replace-InCell 2 6 3 "binary(128)" "binary(16)"

# jelo_masa: integer -> int
Replace-InCell 5 7 3 "integer" "int"

# TIPJELA table: tip + jela + _id  /  tip + jela + _naziv runs merge together.
Replace-InCell 6 2 2 "tipjela_id" "tipjela_id"
Replace-InCell 6 3 2 "tipjela_naziv" "tipjela_naziv"

# UKUS table: ukus + _id  /  ukus + _naziv runs merge together.
Replace-InCell 7 2 2 "ukus_id" "ukus_id"
Replace-InCell 7 3 2 "ukus_naziv" "ukus_naziv"

# DIJETA table: dijeta + _id  /  dijeta + _naziv runs merge together.
Replace-InCell 8 2 2 "dijeta_id" "dijeta_id"
Replace-InCell 8 3 2 "dijeta_naziv" "dijeta_naziv"

# por_br_osoba: integer -> int
Replace-InCell 9 6 3 "integer" "int"

# stavka_kol: integer -> int
Replace-InCell 10 4 3 "integer" "int"
